$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per data refresh
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Update the view: scroll position and selection
$ws.Range("C19").Select()
